$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns per refreshed crypto data feed
$ws.Range("D2").Value = "22.463.09"
$ws.Range("E2").Value = "  +0.32%  "
$ws.Range("D3").Value = "1.573.15"
$ws.Range("E3").Value = "  +0.01%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("D5").Value = "'1.002"
$ws.Range("D6").Value = "'291.35"
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("D7").Value = "'0.3749"
$ws.Range("E7").Value = "  -0.33%  "
$ws.Range("D8").Value = "'49.94"
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("D9").Value = "'0.3406"
$ws.Range("E9").Value = "  -0.25%  "
$ws.Range("D10").Value = "'1.149"
$ws.Range("E10").Value = "  -1.19%  "
$ws.Range("D11").Value = "'0.07586"
$ws.Range("E11").Value = "  -1.06%  "
$ws.Range("E12").Value = "  -0.15%  "
$ws.Range("E13").Value = "  +0.61%  "
$ws.Range("D14").Value = "'5.994"
$ws.Range("E14").Value = "  +0.39%  "
$ws.Range("D15").Value = "'6.962"
$ws.Range("E15").Value = "  +0.58%  "
$ws.Range("D16").Value = "1.572.97"
$ws.Range("E16").Value = "  -0.19%  "
$ws.Range("D17").Value = "'0.00001125"
$ws.Range("E17").Value = "  -0.99%  "
$ws.Range("D18").Value = "'91.21"
$ws.Range("E18").Value = "  +0.70%  "
$ws.Range("D19").Value = "'0.06737"
$ws.Range("E19").Value = "  -0.37%  "
$ws.Range("E20").Value = "  -0.19%  "
$ws.Range("D21").Value = "'6.282"
$ws.Range("E21").Value = "  +0.83%  "
$ws.Range("D22").Value = "'16.44"
$ws.Range("E22").Value = "  -1.75%  "
$ws.Range("E23").Value = "  +1.48%  "
$ws.Range("D24").Value = "22.470.07"
$ws.Range("E24").Value = "  +0.37%  "
$ws.Range("D25").Value = "'2.327"
$ws.Range("E25").Value = "  -4.35%  "
$ws.Range("D26").Value = "'2.603"
$ws.Range("E26").Value = "  -5.44%  "
$ws.Range("E27").Value = "  -0.48%  "
$ws.Range("D28").Value = "'148.87"
$ws.Range("E28").Value = "  +2.34%  "
$ws.Range("D29").Value = "'5.005"
$ws.Range("E29").Value = "  -1.39%  "
$ws.Range("D30").Value = "'126.08"
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("D31").Value = "1.748.96"
$ws.Range("D32").Value = "'1.043"
$ws.Range("E32").Value = "  +2.91%  "
$ws.Range("D33").Value = "'6.151"
$ws.Range("E33").Value = "  -0.76%  "
$ws.Range("E34").Value = "  -2.20%  "
$ws.Range("D35").Value = "'9.900"
$ws.Range("E35").Value = "  -2.36%  "
$ws.Range("D36").Value = "'0.08454"
$ws.Range("E36").Value = "  -1.31%  "
$ws.Range("D37").Value = "'1.387"
$ws.Range("E37").Value = "  +4.09%  "
$ws.Range("D38").Value = "'0.02469"
$ws.Range("E38").Value = "  -3.55%  "
$ws.Range("D39").Value = "'0.2299"
$ws.Range("E39").Value = "  -0.95%  "
$ws.Range("D40").Value = "'0.06555"
$ws.Range("E40").Value = "  +0.24%  "
$ws.Range("D41").Value = "'5.495"
$ws.Range("E41").Value = "  +0.50%  "
$ws.Range("D42").Value = "'11.41"
$ws.Range("E42").Value = "  -1.72%  "
$ws.Range("D43").Value = "'0.6299"
$ws.Range("E43").Value = "  -2.65%  "
$ws.Range("E46").Value = "  +0.67%  "
$ws.Range("D47").Value = "'0.5885"
$ws.Range("E47").Value = "  -2.58%  "
$ws.Range("D48").Value = "'2.097"
$ws.Range("E48").Value = "  +0.12%  "
$ws.Range("D49").Value = "'130.21"
$ws.Range("E49").Value = "  +3.36%  "
$ws.Range("D50").Value = "'1.232"
$ws.Range("E50").Value = "  -5.68%  "
$ws.Range("D51").Value = "'0.07339"
$ws.Range("E51").Value = "  -0.01%  "

# Row 44/45: Frax and EnergySwap swap ranking positions
$ws.Range("B44").Value = "Frax"
$ws.Range("C44").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D44").Value = "'1.001"
$ws.Range("E44").Value = "  -0.14%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'14.03"
$ws.Range("E45").Value = "  -0.22%  "
